$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 11 new rows right after the current last row (78), pushing nothing
# below since row 78 is currently the last used row. Excel auto-copies the
# row-78 number formatting (date style) down into the new rows, which is
# exactly what we want for column A.
$ws.Rows("79:89").Insert()

# New daily demographic breakdown rows for 2020-04-15 (serial 43936),
# continuing the RACE / ETHNICITY / GENDER blocks already present for each
# prior date.
$data = @(
  @(43936, "RACE",      "White",                      2972),
  @(43936, "RACE",      "Black or African American",  1218),
  @(43936, "RACE",      "Other/Two or More Races",     350),
  @(43936, "RACE",      "Asian",                         88),
  @(43936, "RACE",      "Pending",                     1451),
  @(43936, "ETHNICITY", "Not Hispanic or Latino",      3873),
  @(43936, "ETHNICITY", "Hispanic",                     346),
  @(43936, "ETHNICITY", "Pending",                     1860),
  @(43936, "GENDER",    "Female",                      3109),
  @(43936, "GENDER",    "Male",                        2809),
  @(43936, "GENDER",    "Pending",                      161)
)

$r = 79
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Formula = '=IF(C' + $r + '="Pending","",D' + $r + '/SUMIFS(D:D,A:A,A' + $r + ',B:B,B' + $r + ',C:C,"<>Pending"))'
    $r = $r + 1
}

# Match the refreshed view state: scrolled further down and the new last
# cell (A89) selected.
$win = $excel.ActiveWindow
$win.ScrollRow = 55
$win.ScrollColumn = 1
$null = $ws.Range("A89").Select()
